$d = $word.ActiveDocument

# --- First paragraph: the hidden "**ID__..." topic-id marker ---
$p1 = $d.Paragraphs(1)

# 1) Give the paragraph a (borderless) paragraph border - i.e. just the
#    5pt spacing on all four sides, with no visible rule - matching the
#    "See <hyperlink>." paragraph later in the document.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5

# 2) Increase the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.LeftIndent = 11.25

# 3) The paragraph currently holds two runs: the "**ID__...__ID**" token
#    and a trailing run containing a single space. Drop the space-only
#    run entirely and retarget the topic id embedded in the first run.
$paraRange = $p1.Range
$paraEnd = $paraRange.End

# Paragraph.Range always extends through its own paragraph mark, so the
# mark sits at (paraEnd - 1) and the character immediately preceding it
# is the last "real" content character - here, the trailing space run.
$trailingCharRange = $d.Range($paraEnd - 2, $paraEnd - 1)
if ($trailingCharRange.Text -eq " ") {
    $trailingCharRange.Delete()
}

$paraEnd = $p1.Range.End
$textRange = $d.Range($p1.Range.Start, $paraEnd - 1)
$textRange.Text = "**ID__AFFARS_5316_603_3__ID**"
